$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Normalize Spanish connector words (de/del/la/las/el/los/y) to capitalized
# form in state/municipality names

$ws.Range("B7").Value = "Pabellón De Arteaga"
$ws.Range("B8").Value = "Rincón De Romos"
$ws.Range("B9").Value = "San Francisco De Los Romo"
$ws.Range("B13").Value = "Playas De Rosarito"
$ws.Range("B27").Value = "Amatenango De La Frontera"
$ws.Range("B29").Value = "Comitán De Domínguez"
$ws.Range("B43").Value = "Ocozocoautla De Espinosa"
$ws.Range("B47").Value = "San Cristóbal De Las Casas"
$ws.Range("B73").Value = "Guadalupe Y Calvo"
$ws.Range("B76").Value = "Hidalgo Del Parral"
$ws.Range("B88").Value = "San Francisco Del Oro"
$ws.Range("B104").Value = "San Juan De Sabinas"
$ws.Range("A117").Value = "Ciudad De México"
$ws.Range("B121").Value = "Cuajimalpa De Morelos"
$ws.Range("B151").Value = "San Juan Del Río"
$ws.Range("A157").Value = "Estado De México"
$ws.Range("B157").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B160").Value = "Almoloya De Alquisiras"
$ws.Range("B161").Value = "Almoloya De Juárez"
$ws.Range("B169").Value = "Coacalco De Berriozábal"
$ws.Range("B174").Value = "Ecatepec De Morelos"
$ws.Range("B177").Value = "Ixtapan De La Sal"
$ws.Range("B184").Value = "Naucalpan De Juárez"
$ws.Range("B190").Value = "San Felipe Del Progreso"
$ws.Range("B199").Value = "Tenango Del Valle"
$ws.Range("B203").Value = "Tlalnepantla De Baz"
$ws.Range("B207").Value = "Villa Del Carbón"
$ws.Range("B215").Value = "San Miguel De Allende"
$ws.Range("B216").Value = "Apaseo El Grande"
$ws.Range("B222").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B235").Value = "San Diego De La Unión"
$ws.Range("B237").Value = "San Francisco Del Rincón"
$ws.Range("B239").Value = "San Luis De La Paz"
$ws.Range("B241").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B242").Value = "Silao De La Victoria"
$ws.Range("B246").Value = "Valle De Santiago"
$ws.Range("B251").Value = "Acapulco De Juárez"
$ws.Range("B253").Value = "Ajuchitlán Del Progreso"
$ws.Range("B254").Value = "Alcozauca De Guerrero"
$ws.Range("B258").Value = "Atenango Del Río"
$ws.Range("B259").Value = "Atoyac De Álvarez"
$ws.Range("B260").Value = "Ayutla De Los Libres"
$ws.Range("B262").Value = "Buenavista De Cuéllar"
$ws.Range("B263").Value = "Chilapa De Álvarez"
$ws.Range("B264").Value = "Chilpancingo De Los Bravo"
$ws.Range("B268").Value = "Coyuca De Benítez"
$ws.Range("B269").Value = "Coyuca De Catalán"
$ws.Range("B271").Value = "Cuetzala Del Progreso"
$ws.Range("B272").Value = "Cutzamala De Pinzón"
$ws.Range("B277").Value = "Huitzuco De Los Figueroa"
$ws.Range("B278").Value = "Iguala De La Independencia"
$ws.Range("B279").Value = "Zihuatanejo De Azueta"
$ws.Range("B281").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B293").Value = "Taxco De Alarcón"
$ws.Range("B295").Value = "Técpan De Galeana"
$ws.Range("B297").Value = "Tepecoacuilco De Trujano"
$ws.Range("B299").Value = "Tixtla De Guerrero"
$ws.Range("B302").Value = "Tlapa De Comonfort"
$ws.Range("B309").Value = "Atotonilco De Tula"
$ws.Range("B310").Value = "Atotonilco El Grande"
$ws.Range("B312").Value = "Cuautepec De Hinojosa"
$ws.Range("B314").Value = "Huasca De Ocampo"
$ws.Range("B318").Value = "Jacala De Ledezma"
$ws.Range("B322").Value = "Mixquiahuala De Juárez"
$ws.Range("B323").Value = "Molango De Escamilla"
$ws.Range("B324").Value = "Nopala De Villagrán"
$ws.Range("B325").Value = "Omitlán De Juárez"
$ws.Range("B326").Value = "Pachuca De Soto"
$ws.Range("B328").Value = "Progreso De Obregón"
$ws.Range("B332").Value = "Santiago De Anaya"
$ws.Range("B333").Value = "Santiago Tulantepec De Lugo Guerrero"
$ws.Range("B335").Value = "Tenango De Doria"
$ws.Range("B336").Value = "Tepehuacán De Guerrero"
$ws.Range("B337").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B338").Value = "Tezontepec De Aldama"
$ws.Range("B345").Value = "Tula De Allende"
$ws.Range("B346").Value = "Tulancingo De Bravo"
$ws.Range("B348").Value = "Zacualtipán De Ángeles"
$ws.Range("B352").Value = "Ahualulco De Mercado"
$ws.Range("B360").Value = "Encarnación De Díaz"
$ws.Range("B366").Value = "Jilotlán De Los Dolores"
$ws.Range("B371").Value = "Lagos De Moreno"
$ws.Range("B374").Value = "Ojuelos De Jalisco"
$ws.Range("B376").Value = "San Juan De Los Lagos"
$ws.Range("B378").Value = "San Miguel El Alto"
$ws.Range("B380").Value = "Tamazula De Gordiano"
$ws.Range("B384").Value = "Tepatitlán De Morelos"
$ws.Range("B386").Value = "Tlajomulco De Zúñiga"
$ws.Range("B393").Value = "Unión De San Antonio"
$ws.Range("B394").Value = "Unión De Tula"
$ws.Range("B398").Value = "Zacoalco De Torres"
$ws.Range("B400").Value = "Zapotlán El Grande"
$ws.Range("B419").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B465").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B489").Value = "Puente De Ixtla"
$ws.Range("B494").Value = "Tlaltizapán De Zapata"
$ws.Range("B498").Value = "Zacualpan De Amilpas"
$ws.Range("B501").Value = "Amatlán De Cañas"
$ws.Range("B515").Value = "Mier Y Noriega"
$ws.Range("B519").Value = "San Nicolás De Los Garza"
$ws.Range("B523").Value = "Ayoquezco De Aldama"
$ws.Range("B524").Value = "Coicoyán De Las Flores"
$ws.Range("B526").Value = "Cuilápam De Guerrero"
$ws.Range("B527").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B528").Value = "Heroica Ciudad De Huajuapan De León"
$ws.Range("B529").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B534").Value = "Mariscala De Juárez"
$ws.Range("B536").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B539").Value = "Oaxaca De Juárez"
$ws.Range("B540").Value = "Ocotlán De Morelos"
$ws.Range("B541").Value = "Pinotepa De Don Luis"
$ws.Range("B542").Value = "Putla Villa De Guerrero"
$ws.Range("B545").Value = "San Antonino El Alto"
$ws.Range("B555").Value = "San Juan Bautista Lo De Soto"
$ws.Range("B572").Value = "San Pablo Villa De Mitla"
$ws.Range("B573").Value = "San Pedro El Alto"
$ws.Range("B598").Value = "Santo Domingo De Morelos"
$ws.Range("B603").Value = "Tataltepec De Valdés"
$ws.Range("B604").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B605").Value = "Tlacolula De Matamoros"
$ws.Range("B606").Value = "Villa De Tamazulápam Del Progreso"
$ws.Range("B607").Value = "Villa De Tututepec"
$ws.Range("B608").Value = "Villa Sola De Vega"
$ws.Range("B609").Value = "Zimatlán De Álvarez"
$ws.Range("B615").Value = "Chalchicomula De Sesma"
$ws.Range("B623").Value = "Cuayuca De Andrade"
$ws.Range("B629").Value = "Izúcar De Matamoros"
$ws.Range("B634").Value = "Palmar De Bravo"
$ws.Range("B640").Value = "San Nicolás De Los Ranchos"
$ws.Range("B641").Value = "San Salvador El Seco"
$ws.Range("B644").Value = "Tecali De Herrera"
$ws.Range("B650").Value = "Tepexi De Rodríguez"
$ws.Range("B651").Value = "Tetela De Ocampo"
$ws.Range("B653").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B661").Value = "Amealco De Bonfil"
$ws.Range("B663").Value = "Cadereyta De Montes"
$ws.Range("B667").Value = "Jalpan De Serra"
$ws.Range("B669").Value = "Pinal De Amoles"
$ws.Range("B671").Value = "San Juan Del Río"
$ws.Range("B684").Value = "Ciudad Del Maíz"
$ws.Range("B690").Value = "Mexquitic De Carmona"
$ws.Range("B695").Value = "San Ciro De Acosta"
$ws.Range("B703").Value = "Villa De Arriaga"
$ws.Range("B704").Value = "Villa De La Paz"
$ws.Range("B705").Value = "Villa De Ramos"
$ws.Range("B728").Value = "Nacozari De García"
$ws.Range("B764").Value = "Acuamanala De Miguel Hidalgo"
$ws.Range("B765").Value = "Apetatitlán De Antonio Carvajal"
$ws.Range("B770").Value = "Contla De Juan Cuamatzi"
$ws.Range("B772").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B775").Value = "Papalotla De Xicohténcatl"
$ws.Range("B777").Value = "San Pablo Del Monte"
$ws.Range("B793").Value = "Amatlán De Los Reyes"
$ws.Range("B798").Value = "Castillo De Teayo"
$ws.Range("B806").Value = "Cosamaloapan De Carpio"
$ws.Range("B814").Value = "Ignacio De La Llave"
$ws.Range("B823").Value = "Lerdo De Tejada"
$ws.Range("B825").Value = "Martínez De La Torre"
$ws.Range("B827").Value = "Nanchital De Lázaro Cárdenas Del Río"
$ws.Range("B832").Value = "Ozuluama De Mascareñas"
$ws.Range("B836").Value = "Poza Rica De Hidalgo"
$ws.Range("B841").Value = "Soledad De Doblado"
$ws.Range("B857").Value = "Vega De Alatorre"
$ws.Range("B870").Value = "Concepción Del Oro"
$ws.Range("B885").Value = "Noria De Ángeles"
$ws.Range("B892").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B894").Value = "Villa De Cos"

# Remove trailing metadata/footnote rows (901-905)
$ws.Rows("901:905").Delete()
